# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns,
# and swap the Arweave/OKB rows (39/40) per upstream ranking change.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.920.14'
$ws.Range("E2").Value = '  +5.70%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.978.38'
$ws.Range("E3").Value = '  +3.09%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '581.37'
$ws.Range("E5").Value = '  +2.68%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '153.35'
$ws.Range("E6").Value = '  +7.28%  '

$ws.Range("E7").Value = '  -0.18%  '

$ws.Range("E8").Value = '  +1.67%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.974.29'
$ws.Range("E9").Value = '  +2.93%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.99'
$ws.Range("E10").Value = '  +2.75%  '

$ws.Range("E11").Value = '  +2.26%  '

$ws.Range("E12").Value = '  +3.27%  '

$ws.Range("E13").Value = '  +1.11%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.92'
$ws.Range("E14").Value = '  +6.50%  '

$ws.Range("E15").Value = '  +0.93%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.877.86'
$ws.Range("E16").Value = '  +5.62%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.472.42'
$ws.Range("E17").Value = '  +3.10%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.90'
$ws.Range("E18").Value = '  +5.40%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.982.90'
$ws.Range("E19").Value = '  +3.36%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '448.47'
$ws.Range("E20").Value = '  +4.01%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.67'
$ws.Range("E21").Value = '  +4.55%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.680'
$ws.Range("E22").Value = '  +4.15%  '

$ws.Range("E23").Value = '  +6.91%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '80.95'
$ws.Range("E24").Value = '  +2.41%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.36'
$ws.Range("E25").Value = '  +4.72%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.69'
$ws.Range("E26").Value = '  +7.23%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.20'
$ws.Range("E27").Value = '  +10.01%  '

$ws.Range("E28").Value = '  -0.03%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.40'
$ws.Range("E29").Value = '  +17.09%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.86'
$ws.Range("E30").Value = '  +12.58%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.59'
$ws.Range("E31").Value = '  +3.52%  '

$ws.Range("E32").Value = '  -1.27%  '

$ws.Range("E33").Value = '  +4.41%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '26.73'
$ws.Range("E34").Value = '  +4.79%  '

$ws.Range("E35").Value = '  -0.16%  '

$ws.Range("E36").Value = '  +2.67%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.71'
$ws.Range("E37").Value = '  +6.07%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.08'
$ws.Range("E38").Value = '  +7.95%  '

$ws.Range("B39").Value = 'Arweave'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '44.86'
$ws.Range("E39").Value = '  +14.14%  '

$ws.Range("B40").Value = 'OKB'
$ws.Range("C40").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '48.87'
$ws.Range("E40").Value = '  +0.13%  '

$ws.Range("E41").Value = '  +2.85%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.300'
$ws.Range("E42").Value = '  +12.50%  '

$ws.Range("E44").Value = '  +2.04%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '382.96'
$ws.Range("E45").Value = '  +13.34%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.768.31'
$ws.Range("E46").Value = '  +3.16%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0349'
$ws.Range("E47").Value = '  +4.38%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '134.40'
$ws.Range("E48").Value = '  +0.46%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '23.14'
$ws.Range("E50").Value = '  +8.16%  '

$ws.Range("E51").Value = '  +2.43%  '
